$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1957.6666
$ws.Range("I40").Value = 1800
$ws.Range("J40").Value = 1989.2
$ws.Range("K40").Value = 1800
$ws.Range("L40").Value = 1989.2
$ws.Range("M40").Value = -1625
$ws.Range("N40").Value = -2339.2
$ws.Range("H53").Value = 166667800
$ws.Range("I53").Value = 1106.6666
$ws.Range("K53").Value = 1106.6666
$ws.Range("M53").Value = -469.6666
$ws.Range("H64").Value = 4560
$ws.Range("I64").Value = 4266.6665
$ws.Range("K64").Value = 4266.6665
$ws.Range("M64").Value = -4018.6665
$ws.Range("H67").Value = 4560
$ws.Range("I67").Value = 4266.6665
$ws.Range("K67").Value = 4266.6665
$ws.Range("M67").Value = -3408.6665
$ws.Range("H74").Value = 6492.067
$ws.Range("I74").Value = 6798.5386
$ws.Range("K74").Value = 6798.5386
$ws.Range("M74").Value = -5862.5386
$ws.Range("H77").Value = 6492.067
$ws.Range("I77").Value = 6798.5386
$ws.Range("K77").Value = 33992.693
$ws.Range("M77").Value = -29312.693
$ws.Range("H101").Value = 1520.2858
$ws.Range("I101").Value = 1558.9
$ws.Range("K101").Value = 4676.700000000001
$ws.Range("M101").Value = -3054.700000000001
$ws.Range("H112").Value = 5173.5264
$ws.Range("J112").Value = 5664.5293
$ws.Range("L112").Value = 16993.5879
$ws.Range("N112").Value = -19209.5879
$ws.Range("H118").Value = 3190.818
$ws.Range("J118").Value = 9916.666999999999
$ws.Range("L118").Value = 29750.001
$ws.Range("N118").Value = -33064.001
$ws.Range("H129").Value = 6038.1333
$ws.Range("I129").Value = 1057.2
$ws.Range("J129").Value = 16000
$ws.Range("K129").Value = 3171.6
$ws.Range("L129").Value = 48000
$ws.Range("M129").Value = 1828.4
$ws.Range("N129").Value = -58000
$ws.Range("H138").Value = 6749.2188
$ws.Range("J138").Value = 5838.5557
$ws.Range("L138").Value = 17515.6671
$ws.Range("N138").Value = -27795.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1363.3636
$ws.Range("I2").Value = 1028.4286
$ws.Range("K2").Value = 1028.4286
$ws.Range("M2").Value = -915.4286
$ws.Range("H45").Value = 81479.69500000001
$ws.Range("I45").Value = 144905.58
$ws.Range("K45").Value = 144905.58
$ws.Range("M45").Value = -144528.58
$ws.Range("H61").Value = 1454014.9
$ws.Range("I61").Value = 44100.58
$ws.Range("J61").Value = 3745125.5
$ws.Range("K61").Value = 44100.58
$ws.Range("L61").Value = 3745125.5
$ws.Range("M61").Value = -43888.58
$ws.Range("N61").Value = -3745549.5
$ws.Range("H80").Value = 51999.5
$ws.Range("J80").Value = 51999.5
$ws.Range("L80").Value = 51999.5
$ws.Range("N80").Value = -53995.5
$ws.Range("H82").Value = 29000
$ws.Range("J82").Value = 29000
$ws.Range("L82").Value = 29000
$ws.Range("N82").Value = -29722
$ws.Range("H83").Value = 51999.5
$ws.Range("J83").Value = 51999.5
$ws.Range("L83").Value = 155998.5
$ws.Range("N83").Value = -165982.5
$ws.Range("H85").Value = 29000
$ws.Range("J85").Value = 29000
$ws.Range("L85").Value = 29000
$ws.Range("N85").Value = -31496
$ws.Range("H88").Value = 1523.619
$ws.Range("I88").Value = 1900
$ws.Range("J88").Value = 1484
$ws.Range("K88").Value = 1900
$ws.Range("L88").Value = 1484
$ws.Range("M88").Value = -1494
$ws.Range("N88").Value = -2296
$ws.Range("H91").Value = 1523.619
$ws.Range("I91").Value = 1900
$ws.Range("J91").Value = 1484
$ws.Range("K91").Value = 1900
$ws.Range("L91").Value = 1484
$ws.Range("M91").Value = -496
$ws.Range("N91").Value = -4292
$ws.Range("H116").Value = 1363.3636
$ws.Range("I116").Value = 1028.4286
$ws.Range("K116").Value = 1028.4286
$ws.Range("M116").Value = 1265.5714
$ws.Range("H136").Value = 1454014.9
$ws.Range("I136").Value = 44100.58
$ws.Range("J136").Value = 3745125.5
$ws.Range("K136").Value = 132301.74
$ws.Range("L136").Value = 11235376.5
$ws.Range("M136").Value = -129751.74
$ws.Range("N136").Value = -11240476.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1363.3636
$ws.Range("I3").Value = 1028.4286
$ws.Range("K3").Value = 1028.4286
$ws.Range("M3").Value = -914.4286
$ws.Range("H61").Value = 19000
$ws.Range("J61").Value = 19000
$ws.Range("L61").Value = 19000
$ws.Range("N61").Value = -19626
$ws.Range("H86").Value = 7596.5
$ws.Range("I86").Value = 6824.5
$ws.Range("J86").Value = 7927.357
$ws.Range("K86").Value = 6824.5
$ws.Range("L86").Value = 7927.357
$ws.Range("M86").Value = -5701.5
$ws.Range("N86").Value = -10173.357
$ws.Range("H89").Value = 7596.5
$ws.Range("I89").Value = 6824.5
$ws.Range("J89").Value = 7927.357
$ws.Range("K89").Value = 34122.5
$ws.Range("L89").Value = 39636.785
$ws.Range("M89").Value = -28506.5
$ws.Range("N89").Value = -50868.785
$ws.Range("H134").Value = 18001458
$ws.Range("I134").Value = 1363.9395
$ws.Range("K134").Value = 4091.8185
$ws.Range("M134").Value = -1556.8185

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2466.7844
$ws.Range("J31").Value = 5273.778
$ws.Range("L31").Value = 5273.778
$ws.Range("N31").Value = -5863.778
$ws.Range("H32").Value = 12382.2
$ws.Range("J32").Value = 16955.5
$ws.Range("L32").Value = 16955.5
$ws.Range("N32").Value = -17587.5
$ws.Range("H34").Value = 2466.7844
$ws.Range("J34").Value = 5273.778
$ws.Range("L34").Value = 5273.778
$ws.Range("N34").Value = -5677.778
$ws.Range("H99").Value = 3334581.8
$ws.Range("J99").Value = 1063.3334
$ws.Range("L99").Value = 1063.3334
$ws.Range("N99").Value = -4059.3334
$ws.Range("H126").Value = 3334581.8
$ws.Range("J126").Value = 1063.3334
$ws.Range("L126").Value = 3190.0002
$ws.Range("N126").Value = -8130.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2368
$ws.Range("I68").Value = 1803.25
$ws.Range("K68").Value = 5409.75
$ws.Range("M68").Value = -4598.75
$ws.Range("H71").Value = 2368
$ws.Range("I71").Value = 1803.25
$ws.Range("K71").Value = 16229.25
$ws.Range("M71").Value = -12173.25
$ws.Range("H131").Value = 2600390.8
$ws.Range("I131").Value = 10102608
$ws.Range("J131").Value = 3469.1924
$ws.Range("K131").Value = 30307824
$ws.Range("L131").Value = 10407.5772
$ws.Range("M131").Value = -30302784
$ws.Range("N131").Value = -20487.5772

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1974198.2
$ws.Range("I132").Value = 1553
$ws.Range("K132").Value = 4659
$ws.Range("M132").Value = -2129

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4115.095
$ws.Range("I22").Value = 1149.2
$ws.Range("K22").Value = 1149.2
$ws.Range("M22").Value = -854.2
$ws.Range("H27").Value = 4115.095
$ws.Range("I27").Value = 1149.2
$ws.Range("K27").Value = 1149.2
$ws.Range("M27").Value = -1042.2
$ws.Range("H34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()
$ws.Range("H46").Value = 8449.049999999999
$ws.Range("I46").Value = 35999.332
$ws.Range("K46").Value = 35999.332
$ws.Range("M46").Value = -35811.332
$ws.Range("H55").Value = 435.3871
$ws.Range("I55").Value = 301.82352
$ws.Range("K55").Value = 301.82352
$ws.Range("M55").Value = -128.82352
$ws.Range("H100").Value = 3406.3333
$ws.Range("I100").Value = 3406.3333
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3406.3333
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2865.3333
$ws.Range("N100").ClearContents()
$ws.Range("H136").Value = 52277.25
$ws.Range("J136").Value = 2784.75
$ws.Range("L136").Value = 8354.25
$ws.Range("N136").Value = -13454.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 70000
$ws.Range("J27").Value = 70000
$ws.Range("L27").Value = 70000
$ws.Range("N27").Value = -70138
$ws.Range("H100").Value = 802.25
$ws.Range("I100").Value = 688.75
$ws.Range("K100").Value = 1377.5
$ws.Range("M100").Value = -836.5
$ws.Range("H136").Value = 766.75
$ws.Range("I136").Value = 682.9
$ws.Range("K136").Value = 2048.7
$ws.Range("M136").Value = 501.3000000000002
